$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chip-seq dna-binding proteins")

# Insert 3 new columns before column B (old B:D -> new E:G), shifting
# everything from B onward three columns to the right.
$ws.Columns("B:D").Insert()

# Populate the three new header cells.
$ws.Cells.Item(1, 2).Value = "Experiment Alias"
$ws.Cells.Item(1, 3).Value = "Project"
$ws.Cells.Item(1, 4).Value = "Secondary Project"

Write-Output "done"
